$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "rosters"

# Remove PPG, FPG, TPG columns (old I, K, M) leaving GamesPlayed, TotalPts, TotalFouls, TotalTechs
$ws.Range("I1").EntireColumn.Delete()
$ws.Range("J1").EntireColumn.Delete()
$ws.Range("K1").EntireColumn.Delete()

# Fill in player stats
$ws.Range("E2").Value = 25
$ws.Range("G2").Value = 25
$ws.Range("H2").Value = 712
$ws.Range("I2").Value = 75
$ws.Range("J2").Value = 1
$ws.Range("G3").Value = 19
$ws.Range("H3").Value = 514
$ws.Range("I3").Value = 11
$ws.Range("J3").Value = 3
$ws.Range("E4").Value = 18
$ws.Range("G4").Value = 19
$ws.Range("H4").Value = 477
$ws.Range("I4").Value = 61
$ws.Range("J4").Value = 5
$ws.Range("E5").Value = 20
$ws.Range("G5").Value = 16
$ws.Range("H5").Value = 360
$ws.Range("I5").Value = 43
$ws.Range("J5").Value = 2
$ws.Range("E6").Value = 23
$ws.Range("G6").Value = 25
$ws.Range("H6").Value = 520
$ws.Range("I6").Value = 9
$ws.Range("J6").Value = 1
$ws.Range("E7").Value = 41
$ws.Range("G7").Value = 27
$ws.Range("H7").Value = 555
$ws.Range("I7").Value = 65
$ws.Range("J7").Value = 0
$ws.Range("E8").Value = 34
$ws.Range("G8").Value = 27
$ws.Range("H8").Value = 535
$ws.Range("I8").Value = 55
$ws.Range("J8").Value = 0
$ws.Range("E9").Value = 50
$ws.Range("G9").Value = 28
$ws.Range("H9").Value = 533
$ws.Range("I9").Value = 76
$ws.Range("J9").Value = 1
$ws.Range("E10").Value = 28
$ws.Range("G10").Value = 27
$ws.Range("H10").Value = 505
$ws.Range("I10").Value = 35
$ws.Range("J10").Value = 2
$ws.Range("E11").Value = 68
$ws.Range("G11").Value = 14
$ws.Range("H11").Value = 249
$ws.Range("I11").Value = 68
$ws.Range("J11").Value = 3
$ws.Range("E12").Value = 5
$ws.Range("G12").Value = 22
$ws.Range("H12").Value = 390
$ws.Range("I12").Value = 6
$ws.Range("J12").Value = 4
$ws.Range("G13").Value = 21
$ws.Range("H13").Value = 334
$ws.Range("I13").Value = 12
$ws.Range("J13").Value = 5
$ws.Range("E14").Value = 47
$ws.Range("G14").Value = 23
$ws.Range("H14").Value = 352
$ws.Range("I14").Value = 43
$ws.Range("J14").Value = 6
$ws.Range("E15").Value = 27
$ws.Range("G15").Value = 24
$ws.Range("H15").Value = 356
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 7
$ws.Range("E16").Value = 21
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 199
$ws.Range("I16").Value = 51
$ws.Range("J16").Value = 8
$ws.Range("E17").Value = 1
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 252
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 15
$ws.Range("E18").Value = 52
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = 264
$ws.Range("I18").Value = 71
$ws.Range("J18").Value = 0
$ws.Range("E19").Value = 6
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 339
$ws.Range("I19").Value = 69
$ws.Range("J19").Value = 0
$ws.Range("E20").Value = 39
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 252
$ws.Range("I20").Value = 53
$ws.Range("J20").Value = 0
$ws.Range("E21").Value = 16
$ws.Range("G21").Value = 26
$ws.Range("H21").Value = 299
$ws.Range("I21").Value = 10
$ws.Range("J21").Value = 0
$ws.Range("E22").Value = 31
$ws.Range("G22").Value = 13
$ws.Range("H22").Value = 149
$ws.Range("I22").Value = 15
$ws.Range("J22").Value = 1
$ws.Range("E23").Value = 57
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = 159
$ws.Range("I23").Value = 64
$ws.Range("J23").Value = 2
$ws.Range("E24").Value = 7
$ws.Range("G24").Value = 27
$ws.Range("H24").Value = 306
$ws.Range("I24").Value = 56
$ws.Range("J24").Value = 1
$ws.Range("E25").Value = 24
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 249
$ws.Range("I25").Value = 73
$ws.Range("J25").Value = 24
$ws.Range("E26").Value = 55
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 162
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 1
$ws.Range("E27").Value = 33
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 58
$ws.Range("I27").Value = 65
$ws.Range("J27").Value = 5
$ws.Range("E28").Value = 11
$ws.Range("G28").Value = 25
$ws.Range("H28").Value = 241
$ws.Range("I28").Value = 39
$ws.Range("J28").Value = 2
$ws.Range("E29").Value = 2
$ws.Range("G29").Value = 17
$ws.Range("H29").Value = 142
$ws.Range("I29").Value = 59
$ws.Range("J29").Value = 3
$ws.Range("E30").Value = 39
$ws.Range("G30").Value = 21
$ws.Range("H30").Value = 173
$ws.Range("I30").Value = 26
$ws.Range("J30").Value = 1
$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 19
$ws.Range("H31").Value = 148
$ws.Range("I31").Value = 71
$ws.Range("J31").Value = 0
$ws.Range("E32").Value = 42
$ws.Range("G32").Value = 20
$ws.Range("H32").Value = 154
$ws.Range("I32").Value = 65
$ws.Range("J32").Value = 0
$ws.Range("E33").Value = 10
$ws.Range("G33").Value = 25
$ws.Range("H33").Value = 184
$ws.Range("I33").Value = 28
$ws.Range("J33").Value = 0
$ws.Range("E34").Value = 37
$ws.Range("G34").Value = 16
$ws.Range("H34").Value = 107
$ws.Range("I34").Value = 41
$ws.Range("J34").Value = 3
$ws.Range("E35").Value = 30
$ws.Range("G35").Value = 26
$ws.Range("H35").Value = 170
$ws.Range("I35").Value = 44
$ws.Range("J35").Value = 4
$ws.Range("E36").Value = 55
$ws.Range("G36").Value = 14
$ws.Range("H36").Value = 83
$ws.Range("I36").Value = 73
$ws.Range("J36").Value = 34
$ws.Range("E37").Value = 32
$ws.Range("G37").Value = 23
$ws.Range("H37").Value = 135
$ws.Range("I37").Value = 26
$ws.Range("J37").Value = 5
$ws.Range("E38").Value = 17
$ws.Range("G38").Value = 24
$ws.Range("H38").Value = 137
$ws.Range("I38").Value = 66
$ws.Range("J38").Value = 5
$ws.Range("E39").Value = 8
$ws.Range("G39").Value = 23
$ws.Range("H39").Value = 109
$ws.Range("I39").Value = 79
$ws.Range("J39").Value = 1
$ws.Range("E40").Value = 32
$ws.Range("G40").Value = 12
$ws.Range("H40").Value = 55
$ws.Range("I40").Value = 50
$ws.Range("J40").Value = 5
$ws.Range("E41").Value = 26
$ws.Range("G41").Value = 17
$ws.Range("H41").Value = 76
$ws.Range("I41").Value = 8
$ws.Range("J41").Value = 1
$ws.Range("E42").Value = 14
$ws.Range("G42").Value = 28
$ws.Range("H42").Value = 125
$ws.Range("I42").Value = 45
$ws.Range("J42").Value = 0
$ws.Range("E43").Value = 40
$ws.Range("G43").Value = 22
$ws.Range("H43").Value = 92
$ws.Range("I43").Value = 25
$ws.Range("J43").Value = 4
$ws.Range("E44").Value = 44
$ws.Range("G44").Value = 26
$ws.Range("H44").Value = 99
$ws.Range("I44").Value = 80
$ws.Range("J44").Value = 0
$ws.Range("E45").Value = 4
$ws.Range("G45").Value = 27
$ws.Range("H45").Value = 100
$ws.Range("I45").Value = 76
$ws.Range("J45").Value = 1
$ws.Range("E46").Value = 56
$ws.Range("G46").Value = 24
$ws.Range("H46").Value = 68
$ws.Range("I46").Value = 67
$ws.Range("J46").Value = 7
$ws.Range("E47").Value = 38
$ws.Range("G47").Value = 23
$ws.Range("H47").Value = 64
$ws.Range("I47").Value = 21
$ws.Range("J47").Value = 4
$ws.Range("E48").Value = 46
$ws.Range("G48").Value = 24
$ws.Range("H48").Value = 54
$ws.Range("I48").Value = 44
$ws.Range("J48").Value = 8
$ws.Range("E49").Value = 35
$ws.Range("G49").Value = 22
$ws.Range("H49").Value = 28
$ws.Range("I49").Value = 62
$ws.Range("J49").Value = 6

# Update view state
$ws.Range("F2:F49").Select()
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("E42").Select()

# Column widths for new numeric columns
$ws.Range("G1").EntireColumn.AutoFit()
$ws.Range("J1").EntireColumn.AutoFit()
